# ---------------------------------------------------------------------------
# Re-upload of the "Torre de Expansão" Power Query export: the underlying
# SharePoint list query was refreshed (query (33) -> query (34)), picking up
# six new rows and two tardy "Data Fim" values on already-present rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet / query artifacts -----------------------------------
$ws.Name = "base"

$defName = $wb.Names.Item(1)
$defName.Name = "query__34"

$lo = $ws.ListObjects.Item(1)
$lo.Name = "Tabela_query__34"
$qt = $lo.QueryTable
$qt.Name = "query (34)"

# --- backfill two "Data Fim" cells that had since closed out -----------
$ws.Range("L48").Value = 46006
$ws.Range("L52").Value = 46003

# --- append the six freshly-synced rows (159-164) -----------------------
$ws.Rows(159).Insert()
$ws.Range("A159").Value = 'Larissa'
$ws.Range("B159").Value = 'Reclamação'
$ws.Range("C159").Value = 'REDE'
$ws.Range("D159").Value = 3125831000501
$ws.Range("E159").Value = 'Mendes Combustiveis Ltda'
$ws.Range("F159").Value = 'GRUPO MENDES CATATAU - Já possui maquina da REDE em outros postos, porém não consegue fazer a solicitação das maquinas no Mendes Combustiveis Ltda - 3125831000501.'
$ws.Range("G159").Value = 'NÃO'
$ws.Range("H159").Value = 'Central REDE'
$ws.Range("I159").Value = 'Whatsapp'
$ws.Range("J159").Value = 'Receptivo'
$ws.Range("K159").Value = 45988
$ws.Range("M159").Value = 'G.N. Urbano Fortaleza'
$ws.Range("N159").Value = 'CN'
$ws.Range("Q159").Value = 0
$ws.Range("S159").Value = 'Item'
$ws.Range("T159").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'
$ws.Rows(159).RowHeight = 29

$ws.Rows(160).Insert()
$ws.Range("A160").Value = 'Roberta'
$ws.Range("B160").Value = 'Solicitação'
$ws.Range("C160").Value = 'Dúvida PPI'
$ws.Range("D160").Value = 48346249000180
$ws.Range("E160").Value = 'MD Auto posto'
$ws.Range("F160").Value = 'Joyce revendedora solicita o retorno do CT sobre pedido de cancelamento do PPI'
$ws.Range("G160").Value = 'SIM'
$ws.Range("H160").Value = 'Trade Marketing'
$ws.Range("I160").Value = 'Whatsapp'
$ws.Range("J160").Value = 'Receptivo'
$ws.Range("K160").Value = 46006
$ws.Range("L160").Value = 46006
$ws.Range("M160").Value = 'G.N. Rodovia Goiania'
$ws.Range("N160").Value = 'Revendedor'
$ws.Range("Q160").Value = 0
$ws.Range("S160").Value = 'Item'
$ws.Range("T160").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'

$ws.Rows(161).Insert()
$ws.Range("A161").Value = 'Ricardo'
$ws.Range("B161").Value = 'Reclamação'
$ws.Range("C161").Value = 'Ajuste Lat Long'
$ws.Range("D161").Value = 32249473000113
$ws.Range("E161").Value = 'J A Aires Comercial De Combusti Ltda'
$ws.Range("F161").Value = 'A Revendedora Marilu entrou em contato, pois seus clientes não estavam conseguindo favoritar o posto no app KMV. Caso resolvido!'
$ws.Range("G161").Value = 'SIM'
$ws.Range("H161").Value = 'Torre de Expansão'
$ws.Range("I161").Value = 'Whatsapp'
$ws.Range("J161").Value = 'Receptivo'
$ws.Range("K161").Value = 46006
$ws.Range("L161").Value = 46006
$ws.Range("M161").Value = 'G.N. Urbano Santa Maria'
$ws.Range("N161").Value = 'Revendedor'
$ws.Range("Q161").Value = 0
$ws.Range("S161").Value = 'Item'
$ws.Range("T161").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'
$ws.Rows(161).RowHeight = 29

$ws.Rows(162).Insert()
$ws.Range("A162").Value = 'Larissa'
$ws.Range("B162").Value = 'Dúvida'
$ws.Range("C162").Value = 'B2C - Transação Negada'
$ws.Range("D162").Value = 12090263000118
$ws.Range("E162").Value = 'Auto Posto Rio Cervo Ltda'
$ws.Range("F162").Value = 'Cliente relata que um dos seus clientes não consegue fazer transações com o App, não passou informações do cleinte, as informei o que pode ser e o numero da central'
$ws.Range("G162").Value = 'SIM'
$ws.Range("H162").Value = 'Central Atendimento B2C'
$ws.Range("I162").Value = 'Whatsapp'
$ws.Range("J162").Value = 'Receptivo'
$ws.Range("K162").Value = 46006
$ws.Range("L162").Value = 46006
$ws.Range("M162").Value = 'G.N. Urbano Uberlandia'
$ws.Range("N162").Value = 'Revendedor'
$ws.Range("Q162").Value = 0
$ws.Range("S162").Value = 'Item'
$ws.Range("T162").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'
$ws.Rows(162).RowHeight = 29

$ws.Rows(163).Insert()
$ws.Range("A163").Value = 'Ricardo'
$ws.Range("B163").Value = 'Dúvida'
$ws.Range("C163").Value = 'Dúvidas Acúmulos'
$ws.Range("D163").Value = 91411256004800
$ws.Range("E163").Value = ' Sander Comercio De Combustiveis Ltda'
$ws.Range("F163").Value = 'Reunião por chamada de vídeo (Teams) com o revendedor Luan, a fim de esclarecer dúvidas sobre acúmulos e o programa Ipiranga Top.'
$ws.Range("G163").Value = 'SIM'
$ws.Range("H163").Value = 'Torre de Expansão'
$ws.Range("I163").Value = 'Teams'
$ws.Range("J163").Value = 'Ativo'
$ws.Range("K163").Value = 46006
$ws.Range("L163").Value = 46006
$ws.Range("M163").Value = 'G.N. Urbano Passo Fundo'
$ws.Range("N163").Value = 'Revendedor'
$ws.Range("Q163").Value = 0
$ws.Range("S163").Value = 'Item'
$ws.Range("T163").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'
$ws.Rows(163).RowHeight = 29

$ws.Rows(164).Insert()
$ws.Range("A164").Value = 'Larissa'
$ws.Range("B164").Value = 'Solicitação'
$ws.Range("C164").Value = 'Treinamento Financeiro'
$ws.Range("D164").Value = 45431429000171
$ws.Range("E164").Value = 'Afgm Com Derivados De Petroleo Ltda'
$ws.Range("F164").Value = 'Afgm Com Derivados De Petroleo Ltda - 45431429000171| Solicitação de treinamento a pedido do CN marconio (Contato; Daniela Telefone: 31-99324-0071)'
$ws.Range("G164").Value = 'NÃO'
$ws.Range("H164").Value = 'Torre de Expansão'
$ws.Range("I164").Value = 'E-mail'
$ws.Range("J164").Value = 'Receptivo'
$ws.Range("K164").Value = 46006
$ws.Range("M164").Value = 'G.N. Urbano Belo Horizonte'
$ws.Range("N164").Value = 'CN'
$ws.Range("Q164").Value = 0
$ws.Range("S164").Value = 'Item'
$ws.Range("T164").Value = 'personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos'
$ws.Rows(164).RowHeight = 29

# --- grow the query table / autofilter to cover the new rows -----------
$lo.Resize($ws.Range("A1:T164"))
$defName.RefersTo = "=base!`$A`$1:`$T`$164"
$ws.Range("A1:T164").Select()
